$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "move FOOTER_ID to top of footer" ---
# In the FOOTER row (row 10), the third block used to read:
#   K10:P10 = "OBJECT_ID[23:0]"   Q10:R10 = "FOOTER_ID[7:0](0x55)"
# It now reads (single merged block, wider OBJECT_ID field, FOOTER_ID moved to the
# front of the row, replacing the old duplicate ZERO_PADING[7:0] label at C10):
#   C10:D10 = "FOOTER_ID[7:0](0x55)"   K10:R10 = "OBJECT_ID[31:0]"

# 1. Unmerge the two existing merged ranges that will be combined/replaced
$ws.Range("K10:P10").UnMerge()
$ws.Range("Q10:R10").UnMerge()

# 2. Update the text values
#    - C10 (start of footer row) now shows the FOOTER_ID text (was ZERO_PADING[7:0])
#    - Q10's old FOOTER_ID text is cleared since it has moved to C10
#    - A19 "Last Update" note is refreshed
#    - K10 now shows the widened OBJECT_ID[31:0] label (was OBJECT_ID[23:0])
$ws.Range("C10").Value = "FOOTER_ID[7:0](0x55)"
$ws.Range("Q10").Value = ""
$ws.Range("A19").Value = "Last Update: 2020/07/15"
$ws.Range("K10").Value = "OBJECT_ID[31:0]"

# 3. Re-merge K10:R10 into a single block
$ws.Range("K10:R10").Merge()

# Restore the thin box border around the newly-merged K10:R10 block (left & top edges),
# matching the same box styling used by the rest of the footer row (e.g. E10:J10).
$xlEdgeLeft = 7
$xlEdgeTop = 8
$xlThin = 2
$footerObjIdRange = $ws.Range("K10:R10")
$footerObjIdRange.Borders.Item($xlEdgeLeft).LineStyle = $xlThin
$footerObjIdRange.Borders.Item($xlEdgeLeft).Weight = $xlThin
$footerObjIdRange.Borders.Item($xlEdgeTop).LineStyle = $xlThin
$footerObjIdRange.Borders.Item($xlEdgeTop).Weight = $xlThin

# 4. Update the selection to match (K10:R10 selected, active cell K10)
$ws.Range("K10:R10").Select()

$wb.Save()
